$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I14").Value = "sd"
$ws.Range("J14").Value = "Statement-non-opinion"

# Row 21: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I21").Value = "sv"
$ws.Range("J21").Value = "Statement-opinion"

# Row 22: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"

# Row 50: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I50").Value = "sd"
$ws.Range("J50").Value = "Statement-non-opinion"

# Row 82: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I82").Value = "sd"
$ws.Range("J82").Value = "Statement-non-opinion"

# Row 83: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I83").Value = "sd"
$ws.Range("J83").Value = "Statement-non-opinion"

$wb.Save()
